$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.753.22"
$ws.Range("E2").Value = "  +5.51%  "

$ws.Range("D3").Value = "3.637.51"
$ws.Range("E3").Value = "  +5.47%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "592.23"
$ws.Range("E5").Value = "  +2.03%  "

$ws.Range("D6").Value = "193.85"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("D7").Value = "0.641"
$ws.Range("E7").Value = "  +1.98%  "

$ws.Range("D8").Value = "3.628.84"
$ws.Range("E8").Value = "  +5.48%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "0.182"
$ws.Range("E10").Value = "  +6.58%  "

$ws.Range("D11").Value = "0.671"
$ws.Range("E11").Value = "  +4.20%  "

$ws.Range("D12").Value = "57.53"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "0.0000308"
$ws.Range("E13").Value = "  +11.50%  "

$ws.Range("D14").Value = "9.88"
$ws.Range("E14").Value = "  +4.29%  "

$ws.Range("D15").Value = "4.210.60"
$ws.Range("E15").Value = "  +5.66%  "

$ws.Range("D16").Value = "20.51"
$ws.Range("E16").Value = "  +7.85%  "

$ws.Range("D17").Value = "3.630.07"
$ws.Range("E17").Value = "  +5.52%  "

$ws.Range("D18").Value = "70.646.37"
$ws.Range("E18").Value = "  +5.54%  "

$ws.Range("E19").Value = "  +5.57%  "

$ws.Range("E20").Value = "  +2.74%  "

$ws.Range("E21").Value = "  +3.87%  "

$ws.Range("D22").Value = "488.31"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("D23").Value = "19.19"
$ws.Range("E23").Value = "  +10.15%  "

$ws.Range("D24").Value = "5.21"
$ws.Range("E24").Value = "  -3.00%  "

$ws.Range("D25").Value = "4.48"
$ws.Range("E25").Value = "  +2.86%  "

$ws.Range("D26").Value = "91.07"
$ws.Range("E26").Value = "  +2.40%  "

$ws.Range("D27").Value = "3.16"
$ws.Range("E27").Value = "  +6.78%  "

$ws.Range("D28").Value = "11.43"
$ws.Range("E28").Value = "  +3.92%  "

$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  +5.88%  "

$ws.Range("D30").Value = "7.97"
$ws.Range("E30").Value = "  +8.94%  "

$ws.Range("D31").Value = "32.75"
$ws.Range("E31").Value = "  +5.09%  "

$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +8.64%  "

$ws.Range("D33").Value = "12.25"
$ws.Range("E33").Value = "  +3.93%  "

$ws.Range("D34").Value = "67.01"
$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("D35").Value = "612.11"
$ws.Range("E35").Value = "  +1.59%  "

$ws.Range("D36").Value = "40.36"
$ws.Range("E36").Value = "  +8.77%  "

$ws.Range("D37").Value = "0.0₃0833"
$ws.Range("E37").Value = "  +11.48%  "

$ws.Range("D38").Value = "0.410"
$ws.Range("E38").Value = "  +5.51%  "

$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  +2.64%  "

$ws.Range("D42").Value = "3.314.59"
$ws.Range("E42").Value = "  +3.71%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +11.47%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.19"
$ws.Range("E44").Value = "  +18.21%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "3.17"
$ws.Range("E45").Value = "  +9.43%  "

$ws.Range("D46").Value = "0.0460"
$ws.Range("E46").Value = "  +7.09%  "

$ws.Range("D47").Value = "9.58"
$ws.Range("E47").Value = "  +10.97%  "

$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("E49").Value = "  +2.78%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "3.24"
$ws.Range("E51").Value = "  +0.80%  "
